# "login & register DONE ."
# Insert two new leading columns (Name, Surname) in front of the existing
# Username / Password / Address / tel. columns, fill them in for the
# existing 10 users, then append 4 brand-new user rows (12-15) with all
# six fields (Name, Surname, Username, Password, Address, tel.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room: insert two blank columns before column A. This shifts
#    the old A:D (Username/Password/Address/tel.) to C:F, carrying their
#    values and styles along automatically.
# ---------------------------------------------------------------------
$ws.Columns("A:B").Insert()

# ---------------------------------------------------------------------
# 2. Headers for the two new columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Surname"

# ---------------------------------------------------------------------
# 3. Name / Surname for the 10 pre-existing users (rows 2-11).
# ---------------------------------------------------------------------
$names = @("Alex","Benjamin","Charlotte","David","Emily","Fiona","George","Hannah","Jack","Lily")
$surnames = @("Anderson","Brown","Carter","Davis","Edwards","Foster","Garcia","Johnson","Smith","Taylor")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $surnames[$i]
}

# ---------------------------------------------------------------------
# 4. Normalise the formatting of the data block (C2:F11) to plain
#    wrap-text / vertically-centered cells (matches columns C:E already,
#    drops the stray number-format flag that used to live on the old
#    tel. column).
# ---------------------------------------------------------------------
$ws.Range("C2:F11").Style = "Normal"
$ws.Range("C2:F11").VerticalAlignment = -4108
$ws.Range("C2:F11").WrapText = $true

# Drop the old trailing placeholder cells that have no place in the new
# layout (empty styled cells in column G / the old stray D12 cell).
$ws.Range("G2:G11").Clear()
$ws.Range("F12").Clear()

# ---------------------------------------------------------------------
# 5. Four brand-new user rows.
# ---------------------------------------------------------------------
$newRows = @(
    @("Mathee","R.","MatheelnwZa007","1234love","276 love u 4ever","098-765-5432"),
    @("Potato","J.op","ABC098","mmmmm","29737feoifeio","098-234-1524"),
    @("rrr","eeee","abc123","1234ui","232erer","099-999-9999"),
    @("jjj","IOIO","youAndMe1","0987abc","werwrw","076-244-5678")
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# ---------------------------------------------------------------------
# 6. Column widths for the shifted/renumbered data columns.
#    (Column D/new already lands on width 13 after the insert, so it
#    needs no adjustment.)
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 16.5
$ws.Columns("F").ColumnWidth = 13.8

Write-Output "Name/Surname columns inserted and new users appended."
